# IST price update 2025-12-22 13:52
# Insert a new "latest check" column at B (existing B:Y shifts right to
# C:Z), stamp its header with the newest timestamp, and carry the current
# price snapshot into the new column for every SKU row except the one
# whose check came back empty this round.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; everything from B onward shifts right.
$ws.Columns("B:B").Insert()

# The Insert() doesn't carry the original column width onto the new blank
# column, so restore it explicitly (20.17 COM units == the sheet's
# standard 21-char stored width, matching every other data column).
$ws.Columns("B:B").ColumnWidth = 20.17

# New snapshot timestamp for the freshly inserted column.
$ws.Range("B1").Value2 = "2025-12-22 19:18"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 21) {
        # This SKU's price check came back empty this round.
        $ws.Cells($r, 2).Value2 = ""
    } else {
        # Carry the latest known price (now shifted into column C) into
        # the new column B.
        $ws.Cells($r, 2).Value2 = $ws.Cells($r, 3).Value2
    }
}
